# Regenerate save_data to use K (strikeouts) instead of Strike# column.
# Updates column G ("K") for rows 2-23 with recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 3
    20 = 1
    21 = 1
    22 = 1
    23 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
